# Implementação outros turnos em AG (cont.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "outros turnos" columns (T:X) for several rows, plus a few
# --- observation notes (X) and one corrected instance count (B43). ---

$ws.Range("T40").Value = 12.19
$ws.Range("U40").Value = 533
$ws.Range("V40").Value = 450
$ws.Range("W40").Value = 33
$ws.Range("X40").Value = "OK"

$ws.Range("T41").Value = 11.85
$ws.Range("U41").Value = 452
$ws.Range("V41").Value = 441
$ws.Range("W41").Value = 34
$ws.Range("X41").Value = "OK"

$ws.Range("T42").Value = 13.53
$ws.Range("U42").Value = 571
$ws.Range("V42").Value = 455
$ws.Range("W42").Value = 37
$ws.Range("X42").Value = "OK"

$ws.Range("B43").Value = 744
$ws.Range("X43").Value = "Não saía da turma 16"

$ws.Range("T44").Value = 12.14
$ws.Range("U44").Value = 331
$ws.Range("V44").Value = 477
$ws.Range("W44").Value = 48
$ws.Range("X44").Value = "OK"

$ws.Range("X45").Value = "Não saía da turma 18"

$ws.Range("T46").Value = 4.86
$ws.Range("U46").Value = 140
$ws.Range("V46").Value = 316
$ws.Range("W46").Value = 26
$ws.Range("X46").Value = "OK"

$ws.Range("T47").Value = 16.19
$ws.Range("U47").Value = 575
$ws.Range("V47").Value = 491
$ws.Range("W47").Value = 46
$ws.Range("X47").Value = "OK"

$ws.Range("T48").Value = 20.48
$ws.Range("U48").Value = 669
$ws.Range("V48").Value = 442
$ws.Range("W48").Value = 36
$ws.Range("X48").Value = "OK"

$ws.Range("T49").Value = 19.88
$ws.Range("U49").Value = 786
$ws.Range("V49").Value = 519
$ws.Range("W49").Value = 52
$ws.Range("X49").Value = "OK"

$ws.Range("X50").Value = "Fazer de madrugada"
$ws.Range("X51").Value = "Fazer de madrugada"
$ws.Range("X52").Value = "Não saía da turma 21"
$ws.Range("X53").Value = "Não saía da turma 18"

$ws.Range("T54").Value = 6.5
$ws.Range("U54").Value = 320
$ws.Range("V54").Value = 421
$ws.Range("W54").Value = 33
$ws.Range("X54").Value = "OK"

$ws.Range("X55").Value = "Não sai do 125/150 - turma 13"
$ws.Range("X56").Value = "Não sai do 147/150 - turma 8"

$ws.Range("T57").Value = 16.32
$ws.Range("U57").Value = 359
$ws.Range("V57").Value = 602
$ws.Range("W57").Value = 59
$ws.Range("X57").Value = "OK"

$ws.Range("T58").Value = 38.47
$ws.Range("U58").Value = 672
$ws.Range("V58").Value = 651
$ws.Range("W58").Value = 83
$ws.Range("X58").Value = "OK"

$ws.Range("T59").Value = 63.04
$ws.Range("U59").Value = 1543
$ws.Range("V59").Value = 600
$ws.Range("W59").Value = 46
$ws.Range("X59").Value = "OK"

$ws.Range("T60").Value = 10.57
$ws.Range("U60").Value = 469
$ws.Range("V60").Value = 567
$ws.Range("W60").Value = 48
$ws.Range("X60").Value = "OK"

$ws.Range("X61").Value = "Não sai de 125/150"
$ws.Range("X62").Value = "Não sai do 125/150 turma 8"

$ws.Range("T63").Value = 12.84
$ws.Range("U63").Value = 359
$ws.Range("V63").Value = 597
$ws.Range("W63").Value = 59
$ws.Range("X63").Value = "OK"

$ws.Range("X64").Value = "Não sai do 149/150 na turma 21"

$ws.Range("T65").Value = 9.22
$ws.Range("U65").Value = 219
$ws.Range("V65").Value = 567
$ws.Range("W65").Value = 66
$ws.Range("X65").Value = "OK"

$ws.Range("T66").Value = 33.13
$ws.Range("U66").Value = 1205
$ws.Range("V66").Value = 592
$ws.Range("W66").Value = 43
$ws.Range("X66").Value = "OK"

$ws.Range("T67").Value = 16.56
$ws.Range("U67").Value = 454
$ws.Range("V67").Value = 616
$ws.Range("W67").Value = 57
$ws.Range("X67").Value = "OK"

$ws.Range("T93").Value = 27.89
$ws.Range("U93").Value = 719
$ws.Range("V93").Value = 920
$ws.Range("W93").Value = 79
$ws.Range("X93").Value = "OK"

# --- Move the active selection to P14 (matches the saved view state). ---
$ws.Range("P14").Select()
